# Updates cryptocurrency Price (column D) and Volume(1h) (column E) values
# with refreshed market data, and swaps the TrustWalletToken/TheSandbox
# rows (41/42) to reflect their new relative ranking, per the scheduled
# GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as text
# (these numeric-looking / percent-looking strings must stay literal text,
# matching the original inline-string cell contents) and then restore the
# cell's style back to the workbook default so no stray formatting is left
# behind.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '23.104.99'
Set-TextValue "E2" '  -0.37%  '
Set-TextValue "D3" '1.590.45'
Set-TextValue "E3" '  -1.09%  '
Set-TextValue "D4" '0.9984'
Set-TextValue "E4" '  -0.11%  '
Set-TextValue "D5" '0.9980'
Set-TextValue "E5" '  -0.23%  '
Set-TextValue "D6" '301.53'
Set-TextValue "E6" '  -0.14%  '
Set-TextValue "D7" '0.3769'
Set-TextValue "E7" '  -0.45%  '
Set-TextValue "D8" '0.3613'
Set-TextValue "E8" '  -1.16%  '
Set-TextValue "D9" '50.99'
Set-TextValue "E9" '  +3.50%  '
Set-TextValue "D10" '1.244'
Set-TextValue "E10" '  -1.64%  '
Set-TextValue "D11" '0.9984'
Set-TextValue "E11" '  -0.09%  '
Set-TextValue "E12" '  -0.47%  '
Set-TextValue "D13" '22.33'
Set-TextValue "E13" '  -3.18%  '
Set-TextValue "D14" '6.522'
Set-TextValue "E14" '  -1.39%  '
Set-TextValue "D15" '7.370'
Set-TextValue "E15" '  -0.99%  '
Set-TextValue "E16" '  -1.25%  '
Set-TextValue "D17" '1.591.55'
Set-TextValue "E17" '  -1.14%  '
Set-TextValue "D18" '92.89'
Set-TextValue "E18" '  +1.45%  '
Set-TextValue "D19" '0.06772'
Set-TextValue "E19" '  -0.42%  '
Set-TextValue "E20" '  -2.62%  '
Set-TextValue "D21" '6.446'
Set-TextValue "E21" '  -1.99%  '
Set-TextValue "D22" '0.9980'
Set-TextValue "E22" '  -0.32%  '
Set-TextValue "D23" '12.79'
Set-TextValue "E23" '  -2.33%  '
Set-TextValue "D24" '23.094.17'
Set-TextValue "E24" '  -0.46%  '
Set-TextValue "D25" '2.388'
Set-TextValue "E25" '  +1.57%  '
Set-TextValue "D26" '2.911'
Set-TextValue "E26" '  +1.12%  '
Set-TextValue "D27" '20.94'
Set-TextValue "E27" '  -0.87%  '
Set-TextValue "D28" '148.69'
Set-TextValue "E28" '  -1.16%  '
Set-TextValue "D29" '5.218'
Set-TextValue "E29" '  -1.37%  '
Set-TextValue "D30" '132.75'
Set-TextValue "E30" '  -0.01%  '
Set-TextValue "E31" '  -1.86%  '
Set-TextValue "D32" '6.685'
Set-TextValue "E32" '  -2.43%  '
Set-TextValue "D33" '1.766.38'
Set-TextValue "E33" '  -1.29%  '
Set-TextValue "D34" '0.9556'
Set-TextValue "E34" '  -1.67%  '
Set-TextValue "D35" '0.07464'
Set-TextValue "E35" '  -3.20%  '
Set-TextValue "D36" '10.11'
Set-TextValue "E36" '  -0.55%  '
Set-TextValue "D37" '0.02687'
Set-TextValue "E37" '  -2.95%  '
Set-TextValue "D38" '0.2501'
Set-TextValue "E38" '  -2.49%  '
Set-TextValue "D39" '0.08784'
Set-TextValue "D40" '6.083'
Set-TextValue "E40" '  -2.99%  '
Set-TextValue "B41" 'TheSandbox'
Set-TextValue "C41" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D41" '0.7069'
Set-TextValue "E41" '  -2.02%  '
Set-TextValue "B42" 'TrustWalletToken'
Set-TextValue "C42" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D42" '1.354'
Set-TextValue "E42" '  -2.56%  '
Set-TextValue "D43" '12.19'
Set-TextValue "E43" '  -4.90%  '
Set-TextValue "D44" '15.02'
Set-TextValue "E44" '  -4.35%  '
Set-TextValue "D45" '0.6478'
Set-TextValue "E45" '  -3.13%  '
Set-TextValue "D46" '0.9968'
Set-TextValue "E46" '  -0.23%  '
Set-TextValue "D47" '3.994'
Set-TextValue "E47" '  +0.24%  '
Set-TextValue "D48" '2.279'
Set-TextValue "E48" '  -1.50%  '
Set-TextValue "D49" '131.39'
Set-TextValue "E49" '  +0.39%  '
Set-TextValue "D50" '0.07901'
Set-TextValue "E50" '  -1.38%  '
Set-TextValue "D51" '1.212'
Set-TextValue "E51" '  +2.66%  '
